# Apply the crypto price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Leading apostrophe forces Excel to store the literal text
    # (matching the source file, where these are inline/shared
    # strings) instead of auto-coercing numeric-looking text into
    # a Number cell. Resetting the Style afterwards drops the
    # quote-prefix formatting flag so the cell keeps the workbook
    # default style, just like the original cell.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" "71.937.89"
Set-TextCell "E2" "  -1.21%  "
Set-TextCell "D3" "2.663.86"
Set-TextCell "E3" "  +1.28%  "
Set-TextCell "E4" "  -0.05%  "
Set-TextCell "D5" "598.74"
Set-TextCell "E5" "  -1.13%  "
Set-TextCell "D6" "174.69"
Set-TextCell "E6" "  -2.31%  "
Set-TextCell "E7" "  +0.02%  "
Set-TextCell "D8" "0.523"
Set-TextCell "E8" "  -0.62%  "
Set-TextCell "D9" "2.661.95"
Set-TextCell "E9" "  +1.10%  "
Set-TextCell "E10" "  -1.88%  "
Set-TextCell "E11" "  +2.44%  "
Set-TextCell "E12" "  +0.77%  "
Set-TextCell "E13" "  -1.22%  "
Set-TextCell "D14" "3.151.31"
Set-TextCell "E14" "  +0.42%  "
Set-TextCell "E15" "  -2.31%  "
Set-TextCell "D16" "71.762.46"
Set-TextCell "E16" "  -1.26%  "
Set-TextCell "E17" "  -1.54%  "
Set-TextCell "D18" "2.667.43"
Set-TextCell "E18" "  +1.24%  "
Set-TextCell "E19" "  +6.31%  "
Set-TextCell "E20" "  +2.08%  "
Set-TextCell "D21" "371.36"
Set-TextCell "E21" "  -3.52%  "
Set-TextCell "D22" "4.17"
Set-TextCell "E22" "  -0.32%  "
Set-TextCell "E23" "  +0.75%  "
Set-TextCell "D24" "72.02"
Set-TextCell "E24" "  -1.45%  "
Set-TextCell "D25" "1.00"
Set-TextCell "E25" "  -0.01%  "
Set-TextCell "D26" "4.32"
Set-TextCell "E26" "  -1.25%  "
Set-TextCell "D27" "9.73"
Set-TextCell "E27" "  -1.27%  "
Set-TextCell "D28" "2.800.02"
Set-TextCell "E28" "  +1.16%  "
Set-TextCell "D29" "1.00"
Set-TextCell "E29" "  -0.02%  "
Set-TextCell "E30" "  +1.17%  "
Set-TextCell "D31" "8.05"
Set-TextCell "E31" "  +0.14%  "
Set-TextCell "D32" "500.64"
Set-TextCell "E32" "  -6.13%  "
Set-TextCell "E33" "  -2.13%  "
Set-TextCell "D35" "0.999"
Set-TextCell "E35" "  -0.21%  "
Set-TextCell "D36" "162.96"
Set-TextCell "E36" "  -0.29%  "
Set-TextCell "D37" "19.54"
Set-TextCell "E37" "  +1.09%  "
Set-TextCell "D38" "19.04"
Set-TextCell "E38" "  -0.41%  "
Set-TextCell "B39" "ImmutableX"
Set-TextCell "C39" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D39" "1.38"
Set-TextCell "E39" "  -2.05%  "
Set-TextCell "B40" "Kaspa"
Set-TextCell "C40" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D40" "0.110"
Set-TextCell "E40" "  -1.40%  "
Set-TextCell "E41" "  -3.25%  "
Set-TextCell "E43" "  -1.60%  "
Set-TextCell "E44" "  -2.16%  "
Set-TextCell "E45" "  -0.06%  "
Set-TextCell "B46" "Aave"
Set-TextCell "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D46" "156.10"
Set-TextCell "E46" "  +3.34%  "
Set-TextCell "B47" "OKB"
Set-TextCell "C47" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D47" "39.47"
Set-TextCell "E47" "  -0.41%  "
Set-TextCell "D48" "0.560"
Set-TextCell "E48" "  +3.22%  "
Set-TextCell "E49" "  +1.11%  "
Set-TextCell "E50" "  +2.00%  "
Set-TextCell "E51" "  -1.78%  "
